# Daily auto push 2026-01-25: insert a new reading row for 2026/01/25 (日)
# at row 694, right after the existing 2026/01/25 block, pushing the
# remaining rows (old 694-735) down by one to 695-736.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 694 (shifts 694:735 -> 695:736,
# and extends the used range / dimension to D736 automatically).
$ws.Rows.Item(694).Insert()

# Populate the newly inserted row. The date/weekday columns are plain
# text that happen to look like dates/numbers, so prefix with a single
# quote to force text entry and avoid Excel's automatic date coercion.
$ws.Range("A694").Value = "'2026/01/25"
$ws.Range("B694").Value = "日"
$ws.Range("C694").Value = 13
$ws.Range("D694").Value = 18
